$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regression-output values (commit: "add Crisis and Credit Allocation").
# These replace the previous coefficient / t-stat strings in B2:D4.
# Some of the new values look like plain numbers ("0.17", "-0.01", "-0.09",
# "0.98", "-0.89"); the source file stores every value in this table as
# text (shared string), so force text entry the same way a user would in
# Excel (apostrophe prefix) to keep them as strings instead of numbers.

$ws.Cells.Item(2,2).Value = "'0.17"
$ws.Cells.Item(3,2).Value = "'-0.01"
$ws.Cells.Item(4,2).Value = "'-0.09"

$ws.Cells.Item(2,3).Value = "44.29***"
$ws.Cells.Item(3,3).Value = "2.21***"
$ws.Cells.Item(4,3).Value = "'0.98"

$ws.Cells.Item(2,4).Value = "'-0.89"
$ws.Cells.Item(3,4).Value = "0.46***"
$ws.Cells.Item(4,4).Value = "0.82*"
